# nuevos experimentos no convexos
#
# Updates the numeric/expression values that are stored as TEXT (shared
# strings) across several worksheets. A plain `.Value = "1.23"` assignment
# would be auto-coerced to a real number by Excel (losing the text type /
# shared-string reuse), so instead we write a text formula that evaluates
# to the desired string (="...") and then convert it in place to a literal
# value with PasteSpecial(xlPasteValues). That preserves the text type
# without touching cell styles.

function Set-TextValue {
    param(
        $Cell,
        [string]$Text
    )
    # Escape embedded double quotes for the formula literal.
    $escaped = $Text.Replace('"', '""')
    $Cell.Formula = '="' + $escaped + '"'
    $Cell.Copy()
    $Cell.PasteSpecial(-4163)  # xlPasteValues
}

$wb = $excel.ActiveWorkbook

# NOTE: worksheet lookup by name is case-insensitive in this object model,
# and this workbook has two sheets whose names differ only by case
# ("Vector_bf" vs "Vector_BF"). Index into Worksheets by position instead
# of by name so each sheet is addressed unambiguously (sheet order matches
# the tab order: 1 Funciones_Objetivo, 2 Restricciones_del_lider,
# 3 Restricciones_del_follower, 4 Punto_modificado, 5 Vector_bf,
# 6 Vector_BF, 7 Vector_Alpha).

# --- Sheet "Restricciones_del_follower" ---------------------------------
$ws3 = $wb.Worksheets.Item(3)

Set-TextValue $ws3.Cells.Item(2,1) "4.449999999999999 - 2x_1 + y_1 - y_2"
Set-TextValue $ws3.Cells.Item(2,2) "-1.9499999999999995"
Set-TextValue $ws3.Cells.Item(2,6) "0"

Set-TextValue $ws3.Cells.Item(3,1) "-1.2499999999999993 + x_1 - 3x_2 + y_2"
Set-TextValue $ws3.Cells.Item(3,2) "-0.7500000000000007"
Set-TextValue $ws3.Cells.Item(3,4) "0.7"
Set-TextValue $ws3.Cells.Item(3,5) "0"
Set-TextValue $ws3.Cells.Item(3,6) "0"

Set-TextValue $ws3.Cells.Item(4,1) "-6.43 + x_1 + x_2"
Set-TextValue $ws3.Cells.Item(4,2) "4.1"
Set-TextValue $ws3.Cells.Item(4,4) "0.8"
Set-TextValue $ws3.Cells.Item(4,5) "0"
Set-TextValue $ws3.Cells.Item(4,6) "0"

# --- Sheet "Punto_modificado" -------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

Set-TextValue $ws4.Cells.Item(2,1) "4.5"
Set-TextValue $ws4.Cells.Item(2,2) "1.6"
Set-TextValue $ws4.Cells.Item(2,3) "6.1000000000000005"
Set-TextValue $ws4.Cells.Item(2,4) "1.55"

# --- Sheet "Vector_bf" ----------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

Set-TextValue $ws5.Cells.Item(3,1) "-0.74"

# --- Sheet "Vector_BF" ----------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

Set-TextValue $ws6.Cells.Item(2,1) "2.0"
Set-TextValue $ws6.Cells.Item(3,1) "-1.0"
Set-TextValue $ws6.Cells.Item(5,1) "-0.0"
